$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move "?Casey" out of E10 down to a new entry in E25 renamed to "Casey",
# and add a new name "Fion" in E26.
$ws.Range("E10").ClearContents()
$ws.Range("E25").Value = "Casey"
$ws.Range("E26").Value = "Fion"

# Update the active selection to match the new state of the sheet.
$ws.Range("F13").Select()
